$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.149.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.246.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.84%  "

$ws.Range("E6").Value = "  -1.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.75"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.51%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0955"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.581.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.240.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.042.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.85%  "

$ws.Range("E19").Value = "  -2.47%  "

$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.89%  "

$ws.Range("E26").Value = "  -7.48%  "

$ws.Range("E27").Value = "  -5.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0828"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").Value = "  -5.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("E41").Value = "  -6.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "113.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.22%  "

$ws.Range("E43").Value = "  -6.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.87%  "

$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("E48").Value = "  -3.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.447"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.51%  "

$ws.Range("E51").Value = "  +0.45%  "
